$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "via only the hotel name. (",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "via only the hotel name to see pre-filled data. (",
    2
)
